$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.61 = 26033.71 pesos`n✅ 26033.71 pesos = 6.58 = 969.66 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the N10/O10 and N12/O12 rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 151.189
$ws2.Range("O10").Value = 3936.01
$ws2.Range("N12").Value = 3957.99
$ws2.Range("O12").Value = 147.42
